$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.334.65'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '2.517.78'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = "'522.56"
$ws.Range('E5').Value = '  +0.33%  '
$ws.Range('D6').Value = "'132.90"
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('D8').Value = "'0.559"
$ws.Range('E8').Value = '  +0.28%  '
$ws.Range('D9').Value = '2.515.37'
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('D10').Value = "'0.0976"
$ws.Range('E10').Value = '  -1.14%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = "'5.15"
$ws.Range('E12').Value = '  -2.40%  '
$ws.Range('E13').Value = '  -2.38%  '
$ws.Range('D14').Value = '2.956.57'
$ws.Range('E14').Value = '  +1.60%  '
$ws.Range('D15').Value = '58.398.97'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = "'22.11"
$ws.Range('E16').Value = '  -1.49%  '
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '2.502.41'
$ws.Range('E18').Value = '  +1.10%  '
$ws.Range('E19').Value = '  -0.53%  '
$ws.Range('D20').Value = "'321.85"
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  -0.70%  '
$ws.Range('D22').Value = "'6.15"
$ws.Range('E22').Value = '  +6.82%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').Value = "'64.40"
$ws.Range('E25').Value = '  -1.88%  '
$ws.Range('D26').Value = "'0.993"
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('E27').Value = '  +0.39%  '
$ws.Range('E28').Value = '  -0.11%  '
$ws.Range('D29').Value = '0.0₃0755'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = "'168.14"
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = "'1.20"
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('D33').Value = "'6.31"
$ws.Range('E33').Value = '  +0.36%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').Value = "'0.997"
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = "'18.08"
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('E37').Value = '  -7.02%  '
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('E39').Value = '  -0.45%  '
$ws.Range('D40').Value = "'36.22"
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').Value = "'0.780"
$ws.Range('E41').Value = '  -3.19%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Value = "'279.03"
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = "'3.49"
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('E44').Value = '  -3.20%  '
$ws.Range('E45').Value = '  +1.07%  '
$ws.Range('D46').Value = "'122.68"
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('D47').Value = "'0.0920"
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('D48').Value = "'0.0502"
$ws.Range('E48').Value = '  +2.10%  '
$ws.Range('D49').Value = "'17.72"
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = "'16.90"
$ws.Range('E51').Value = '  -1.15%  '
